# Generate Report for Handback
#
# This script mirrors a localization "handback" run: the two loc target
# files (zh-cn, de-de) have come back in sync with en-US, so the status
# everywhere flips from "Ready for handoff" to "Handed back: in sync with
# en-US", the per-file "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the language sheets get populated,
# new hyperlinks are added on the "Latest Target File" cells, and a couple
# of columns that now hold longer text get widened.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Helper: translate a "desired stored <col width=.../>" value into the
# ColumnWidth the COM layer needs to be handed so that, after its internal
# pixel-snap/padding round trip, the persisted width lands as close as
# possible to the target. (The engine quantizes to 1/6-character pixel
# steps and adds a constant 5/6 character padding on save.)
function Get-ColumnWidthInput([double]$target) {
    $px = [Math]::Round(($target - 5.0 / 6.0) * 6.0)
    return $px / 6.0
}

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) on both data rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = Get-ColumnWidthInput 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = Get-ColumnWidthInput 29.9777047293527

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn and de-de carry the same shape of update,
# just with different handoff/handback file names and handback timestamp.
# ---------------------------------------------------------------------
function Update-LanguageSheet(
    [string]$sheetName,
    [string]$handbackDateTime
) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Capture the existing A2 / A3 hyperlink target URLs so the new
    # "Latest Target File" (column I) hyperlinks reuse the exact same
    # external targets as the "Source File Name" (column A) links.
    $urlRow2 = $null
    $urlRow3 = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq "`$A`$2") { $urlRow2 = $h.Address }
        if ($addr -eq "`$A`$3") { $urlRow3 = $h.Address }
    }

    # Status column
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File / Latest Handback File / Latest Handback DateTime
    $ws.Range("J2").Value = $ws.Range("G2").Value2
    $ws.Range("J3").Value = $ws.Range("G3").Value2
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime

    $name2 = $ws.Range("A2").Value2
    $name3 = $ws.Range("A3").Value2
    $ws.Range("I2").Value = $name2
    $ws.Range("I3").Value = $name3

    if ($urlRow2) {
        $ws.Hyperlinks.Add($ws.Range("I2"), $urlRow2, "", "", $name2) | Out-Null
    }
    if ($urlRow3) {
        $ws.Hyperlinks.Add($ws.Range("I3"), $urlRow3, "", "", $name3) | Out-Null
    }

    # Column widths: Status (C) and the new Target/Handback file columns (I, J)
    $ws.Columns.Item(3).ColumnWidth = Get-ColumnWidthInput 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = Get-ColumnWidthInput 40
    $ws.Columns.Item(10).ColumnWidth = Get-ColumnWidthInput 40
}

Update-LanguageSheet "zh-cn" "2016-08-21 16:54:33"
Update-LanguageSheet "de-de" "2016-08-21 16:54:39"
